$d = $word.ActiveDocument

# Update the date heading (unique text in the document, Find/Replace is safe)
$d.Content.Find.Execute("2024-12-14 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-15 Sunday", 2) | Out-Null

# Update each arithmetic-answer cell in the practice table, addressed by its
# (row, column) position. We assign Cell.Range.Text directly (rather than
# Find/Replace) so that cells sharing identical source text (e.g. the two
# "23+39=62" cells) are each updated to their own distinct target value.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "83-17=66"
$t.Cell(1, 2).Range.Text = "40+26=66"
$t.Cell(1, 3).Range.Text = "76-60=16"
$t.Cell(1, 4).Range.Text = "79-6=73"
$t.Cell(1, 5).Range.Text = "51-26=25"
$t.Cell(2, 1).Range.Text = "84-22=62"
$t.Cell(2, 2).Range.Text = "69+6=75"
$t.Cell(2, 3).Range.Text = "47-5=42"
$t.Cell(2, 4).Range.Text = "35-6=29"
$t.Cell(2, 5).Range.Text = "46+13=59"
$t.Cell(3, 1).Range.Text = "30-10=20"
$t.Cell(3, 2).Range.Text = "31+38=69"
$t.Cell(3, 3).Range.Text = "75+22=97"
$t.Cell(3, 4).Range.Text = "47+1=48"
$t.Cell(3, 5).Range.Text = "90-63=27"
$t.Cell(4, 1).Range.Text = "50+49=99"
$t.Cell(4, 2).Range.Text = "11+78=89"
$t.Cell(4, 3).Range.Text = "50-32=18"
$t.Cell(4, 4).Range.Text = "31+42=73"
$t.Cell(4, 5).Range.Text = "7+53=60"
$t.Cell(5, 1).Range.Text = "2+4=6"
$t.Cell(5, 2).Range.Text = "18+36=54"
$t.Cell(5, 3).Range.Text = "35+1=36"
$t.Cell(5, 4).Range.Text = "25+2=27"
$t.Cell(5, 5).Range.Text = "76+2=78"
$t.Cell(6, 1).Range.Text = "68+9=77"
$t.Cell(6, 2).Range.Text = "24+45=69"
$t.Cell(6, 3).Range.Text = "44+19=63"
$t.Cell(6, 4).Range.Text = "81-72=9"
$t.Cell(6, 5).Range.Text = "69+7=76"
$t.Cell(7, 1).Range.Text = "83-80=3"
$t.Cell(7, 2).Range.Text = "76+3=79"
$t.Cell(7, 3).Range.Text = "90-9=81"
$t.Cell(7, 4).Range.Text = "65-18=47"
$t.Cell(7, 5).Range.Text = "80-13=67"
$t.Cell(8, 1).Range.Text = "96-90=6"
$t.Cell(8, 2).Range.Text = "2+68=70"
$t.Cell(8, 3).Range.Text = "55+15=70"
$t.Cell(8, 4).Range.Text = "65-8=57"
$t.Cell(8, 5).Range.Text = "20+5=25"
$t.Cell(9, 1).Range.Text = "85-49=36"
$t.Cell(9, 2).Range.Text = "71+18=89"
$t.Cell(9, 3).Range.Text = "56+37=93"
$t.Cell(9, 4).Range.Text = "12-0=12"
$t.Cell(9, 5).Range.Text = "89-75=14"
$t.Cell(10, 1).Range.Text = "69-68=1"
$t.Cell(10, 2).Range.Text = "13+62=75"
$t.Cell(10, 3).Range.Text = "90+2=92"
$t.Cell(10, 4).Range.Text = "23-6=17"
$t.Cell(10, 5).Range.Text = "2+57=59"
$t.Cell(11, 1).Range.Text = "68+23=91"
$t.Cell(11, 2).Range.Text = "9+13=22"
$t.Cell(11, 3).Range.Text = "63+10=73"
$t.Cell(11, 4).Range.Text = "76+10=86"
$t.Cell(11, 5).Range.Text = "82-48=34"
$t.Cell(12, 1).Range.Text = "77+5=82"
$t.Cell(12, 2).Range.Text = "80-0=80"
$t.Cell(12, 3).Range.Text = "94-29=65"
$t.Cell(12, 4).Range.Text = "25+32=57"
$t.Cell(12, 5).Range.Text = "98-82=16"
$t.Cell(13, 1).Range.Text = "31-29=2"
$t.Cell(13, 2).Range.Text = "41+53=94"
$t.Cell(13, 3).Range.Text = "43+30=73"
$t.Cell(13, 4).Range.Text = "55+17=72"
$t.Cell(13, 5).Range.Text = "4+74=78"
$t.Cell(14, 1).Range.Text = "78-74=4"
$t.Cell(14, 2).Range.Text = "69-40=29"
$t.Cell(14, 3).Range.Text = "28+11=39"
$t.Cell(14, 4).Range.Text = "48+41=89"
$t.Cell(14, 5).Range.Text = "35+57=92"
$t.Cell(15, 1).Range.Text = "44+17=61"
$t.Cell(15, 2).Range.Text = "24+39=63"
$t.Cell(15, 3).Range.Text = "78-8=70"
$t.Cell(15, 4).Range.Text = "17+35=52"
$t.Cell(15, 5).Range.Text = "88-25=63"
$t.Cell(16, 1).Range.Text = "23-20=3"
$t.Cell(16, 2).Range.Text = "18+54=72"
$t.Cell(16, 3).Range.Text = "50-8=42"
$t.Cell(16, 4).Range.Text = "78-9=69"
$t.Cell(16, 5).Range.Text = "64-31=33"
$t.Cell(17, 1).Range.Text = "88-31=57"
$t.Cell(17, 2).Range.Text = "77-50=27"
$t.Cell(17, 3).Range.Text = "5+42=47"
$t.Cell(17, 4).Range.Text = "3+52=55"
$t.Cell(17, 5).Range.Text = "94-45=49"
$t.Cell(18, 1).Range.Text = "29+42=71"
$t.Cell(18, 2).Range.Text = "16+71=87"
$t.Cell(18, 3).Range.Text = "84+9=93"
$t.Cell(18, 4).Range.Text = "85-14=71"
$t.Cell(18, 5).Range.Text = "52-11=41"
$t.Cell(19, 1).Range.Text = "36-9=27"
$t.Cell(19, 2).Range.Text = "90-86=4"
$t.Cell(19, 3).Range.Text = "42+30=72"
$t.Cell(19, 4).Range.Text = "90-13=77"
$t.Cell(19, 5).Range.Text = "94-27=67"
$t.Cell(20, 1).Range.Text = "8+62=70"
$t.Cell(20, 2).Range.Text = "32+23=55"
$t.Cell(20, 3).Range.Text = "80-36=44"
$t.Cell(20, 4).Range.Text = "75-40=35"
$t.Cell(20, 5).Range.Text = "95-83=12"
